$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 152 and 153 (QT_SIT_DESVINCULADO and QT_SIT_TRANSFERIDO), shifting
# everything below up by two rows.
$ws.Rows("152:153").Delete()
